# Cruise passes for new PatternUtils strategy
# On the "Results" sheet, a new payment-date row is inserted right after the
# header row (pushing the existing data rows down by one), a new trailing
# row is added at the bottom to keep the same number of "future" rows, and
# the active selection moves to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$ws.Activate()

# Insert a new row above the current row 2 (old first data row), shifting
# rows 2:10 down to 3:11. This also extends the used dimension to B11 and
# auto-creates A11/B11 as part of the shift.
$ws.Rows("2:2").Insert()

# The inserted row is blank/unstyled; copy the formatting that now lives on
# row 3 (the row that used to be row 2) down into the new row 2 so the new
# cells pick up the same number formats/fills as their neighbours.
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)

# Populate the new first data row.
$ws.Range("A2").Value = 40544
$ws.Range("B2").Value = 0

# Move the active cell/selection to B3, matching the post-edit workbook.
$ws.Range("B3").Select()
